# Commit: "Added GetLastNonEmptyRowIndex and GetNonEmptyRowsForWorksheet functions"
#
# Corresponding workbook edits (test fixture used by the new tests):
#   - Switch workbook calculation mode to manual.
#   - Add a new non-empty row (row 4) to "Page3" with cell B4 containing a
#     two-space string - the extra "non-empty" row the new
#     GetLastNonEmptyRowIndex / GetNonEmptyRowsForWorksheet tests look for.
#   - Switch the active sheet from "Page3" to "Page2", and update both
#     sheets' remembered selections.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # Page2
$ws3 = $wb.Worksheets.Item(3)   # Page3

# Workbook goes into manual calculation mode (<calcPr calcMode="manual" .../>).
$excel.Calculation = -4135   # xlCalculationManual

# New trailing "non-empty" row on Page3: B4 = two spaces.
$ws3.Range("B4").Value = "  "

# Page3 keeps its own (now different) selection even though it stops being
# the active sheet.
$ws3.Range("E9").Select() | Out-Null

# Page2 becomes the active / selected sheet, with its own new selection.
$ws2.Activate() | Out-Null
$ws2.Range("C12").Select() | Out-Null
